$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '51.728.03'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.11%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.800.80'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.75%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '356.13'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '109.06'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.556'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("E9").Value = '  +7.13%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '39.87'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.07%  '

# Row 11
$ws.Range("E11").Value = '  +0.29%  '

# Row 12
$ws.Range("E12").Value = '  -0.75%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '19.90'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +2.12%  '

# Row 14
$ws.Range("E14").Value = '  +1.96%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '3.243.33'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.79%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.804.95'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.941'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.35%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '51.690.13'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '

# Row 19
$ws.Range("E19").Value = '  +3.39%  '

# Row 20
$ws.Range("E20").Value = '  +2.92%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.56'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.00%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.0₃0978'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.76%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '70.37'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.20%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '268.24'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.40%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

# Row 26
$ws.Range("E26").Value = '  +0.04%  '

# Row 27
$ws.Range("E27").Value = '  -1.03%  '

# Row 28
$ws.Range("E28").Value = '  -0.10%  '

# Row 29
$ws.Range("E29").Value = '  +1.24%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '37.58'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +8.33%  '

# Row 31
$ws.Range("E31").Value = '  -0.93%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.22'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.98%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '51.93'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.18%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.70'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +11.15%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0443'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.44%  '

# Row 36
$ws.Range("E36").Value = '  +2.39%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '18.81'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.60%  '

# Row 39
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.83%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.14'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.28%  '

# Row 41
$ws.Range("E41").Value = '  +0.97%  '

# Row 42
$ws.Range("E42").Value = '  -2.99%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '22.01'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.53%  '

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.32%  '

# Row 45
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '119.11'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.70%  '

# Row 46
$ws.Range("E46").Value = '  +7.45%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.106.58'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.14%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '3.38'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +3.47%  '

# Row 49
$ws.Range("E49").Value = '  +9.34%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.910'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.73%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '5.41'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -5.73%  '
